# Insert a new two-row block (new weekly price observations) above row 296.
# This shifts all existing rows from 296 downward by two rows, which matches
# the rest of the diff (every later row's content is simply the previous
# row's content, shifted down by two).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A296:A297").EntireRow.Insert()

# Common / constant columns for this sheet's data block.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100102
$producto  = "Cítricos"
$categoriaId = 100102005
$categoria = "Naranja"
$unidad    = "`$/bins (400 kilos)"
$origen    = "Provincia de Limarí"
$kgUnidad  = 400

# New row 296: Fukumoto / Primera
$ws.Cells.Item(296, 1).Value = $mercadoId
$ws.Cells.Item(296, 2).Value = $mercado
$ws.Cells.Item(296, 3).Value = $region
$ws.Cells.Item(296, 4).Value = 44755
$ws.Cells.Item(296, 5).Value = $codreg
$ws.Cells.Item(296, 6).Value = $tipo
$ws.Cells.Item(296, 7).Value = $productoId
$ws.Cells.Item(296, 8).Value = $producto
$ws.Cells.Item(296, 9).Value = $categoriaId
$ws.Cells.Item(296, 10).Value = $categoria
$ws.Cells.Item(296, 11).Value = "Fukumoto"
$ws.Cells.Item(296, 12).Value = "Primera"
$ws.Cells.Item(296, 13).Value = 20
$ws.Cells.Item(296, 14).Value = 115000
$ws.Cells.Item(296, 15).Value = 120000
$ws.Cells.Item(296, 16).Value = 117500
$ws.Cells.Item(296, 17).Value = $unidad
$ws.Cells.Item(296, 18).Value = $origen
$ws.Cells.Item(296, 19).Value = 294
$ws.Cells.Item(296, 20).Value = $kgUnidad

# New row 297: Fukumoto / Segunda
$ws.Cells.Item(297, 1).Value = $mercadoId
$ws.Cells.Item(297, 2).Value = $mercado
$ws.Cells.Item(297, 3).Value = $region
$ws.Cells.Item(297, 4).Value = 44755
$ws.Cells.Item(297, 5).Value = $codreg
$ws.Cells.Item(297, 6).Value = $tipo
$ws.Cells.Item(297, 7).Value = $productoId
$ws.Cells.Item(297, 8).Value = $producto
$ws.Cells.Item(297, 9).Value = $categoriaId
$ws.Cells.Item(297, 10).Value = $categoria
$ws.Cells.Item(297, 11).Value = "Fukumoto"
$ws.Cells.Item(297, 12).Value = "Segunda"
$ws.Cells.Item(297, 13).Value = 20
$ws.Cells.Item(297, 14).Value = 95000
$ws.Cells.Item(297, 15).Value = 100000
$ws.Cells.Item(297, 16).Value = 97500
$ws.Cells.Item(297, 17).Value = $unidad
$ws.Cells.Item(297, 18).Value = $origen
$ws.Cells.Item(297, 19).Value = 244
$ws.Cells.Item(297, 20).Value = $kgUnidad

# Make sure the date column keeps its date number format on the new rows.
$ws.Range("D296:D297").NumberFormat = "YYYY-MM-DD HH:MM:SS"
